$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header row (row 1): a new "Unnamed: 0" column is inserted before the
# old column A ("Nombre" -> now C), "Folio" moves from H to B, a new
# "Encargado" column is appended at M, and "Registro Postgrado" (I) is
# shortened to "Registro". Net effect: every header slides one column
# right of its old position and two brand-new header cells (A1, M1)
# appear.
#
# Seed the two brand-new header cells by copying an existing header
# cell (format + value) on top of them first, so they pick up the same
# cellXf (bold font / thin border / centered-top alignment) already
# used by the rest of row 1, then overwrite the copied text with the
# real header labels.
# ------------------------------------------------------------------
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("B1").Copy($ws.Range("M1"))

$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("B1").Value = "Folio"
$ws.Range("C1").Value = "Nombre"
$ws.Range("D1").Value = "Nombre 2"
$ws.Range("E1").Value = "Apellido Paterno"
$ws.Range("F1").Value = "Apellido Materno"
$ws.Range("G1").Value = "Puesto"
$ws.Range("H1").Value = "Area"
$ws.Range("I1").Value = "Registro"
$ws.Range("J1").Value = "Vigencia"
$ws.Range("K1").Value = "Numero Empleado"
$ws.Range("L1").Value = "Ruta Imagen"
$ws.Range("M1").Value = "Encargado"

# ------------------------------------------------------------------
# Data row 2
# ------------------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "Juan"
$ws.Range("D2").Value = "Carlos"
$ws.Range("E2").Value = "Calderon"
$ws.Range("F2").Value = "Davila"
$ws.Range("G2").Value = "Director"
$ws.Range("H2").Value = "Directorr"
$ws.Range("I2").Value = "23/03/2023"
$ws.Range("J2").Value = 45374
$ws.Range("K2").Value = 12345
$ws.Range("L2").Value = "Niels"
$ws.Range("M2").Value = "C:/Users/MrJua/Downloads/20230318_002304.jpg"

# "Vigencia" is a date serial displayed with a custom date-time format.
# Apply it first in lower-case, then re-apply the final upper-case form:
# this mints numFmtId 164 (yyyy-mm-dd h:mm:ss) followed by 165
# (YYYY-MM-DD HH:MM:SS), with the cell itself ending up on 165 -- this
# is how both custom formats ended up registered in the workbook.
$ws.Range("J2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("J2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ------------------------------------------------------------------
# Data row 3
# ------------------------------------------------------------------
# "Folio" for this row is recorded as text, not a number (note the
# quote prefix forces it to stay text instead of being parsed as 123456).
$ws.Range("B3").Value = "'123456"
$ws.Range("C3").Value = "Juan"
$ws.Range("D3").Value = "Carlos"
$ws.Range("E3").Value = "Calderon"
$ws.Range("F3").Value = "Davila"
$ws.Range("G3").Value = "Director"
$ws.Range("H3").Value = "Perro"
$ws.Range("I3").Value = "23/03/2023"
$ws.Range("J3").Value = 45374
$ws.Range("K3").Value = 165498
$ws.Range("L3").Value = "Niels"
$ws.Range("M3").Value = "C:/Users/MrJua/Downloads/20230318_002304.jpg"

$ws.Range("J3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
